# feat: add 2022-Q1 data
#
# - The existing "总计" sheet (fund-holdings detail previously lived on
#   "2021-Q4", and the running summary lived on "总计") is renamed to
#   "2022-Q1" and repopulated with the 2022-Q1 holdings detail (same
#   column layout as the "2021-Q4" sheet).
# - A brand-new "总计" sheet is appended right after "2022-Q1" (cloned
#   from it, to inherit the original sheet-level formatting) with the
#   running summary table, now covering both 2022-Q1 (first row) and
#   2021-Q4.

$wb = $excel.ActiveWorkbook

$sheetQ4  = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Repurpose the old "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $oldTotal
$q1.Name = "2022-Q1"

# Header row (bold / bordered / centred, matches the "2021-Q4" sheet).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 鹏华创新成长混合A
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'011460"
$q1.Range("C2").Value = "鹏华创新成长混合A"
$q1.Range("D2").Value = "'12.30"
$q1.Range("E2").Value = "'93.30"
$q1.Range("F2").Value = "'4.34"
$q1.Range("G2").Value = "'0.5338"
$q1.Range("H2").Value = 7

# Row 3 - 鹏华创新驱动混合
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'005967"
$q1.Range("C3").Value = "鹏华创新驱动混合"
$q1.Range("D3").Value = "'1.50"
$q1.Range("E3").Value = "'92.83"
$q1.Range("F3").Value = "'4.13"
$q1.Range("G3").Value = "'0.0620"
$q1.Range("H3").Value = 9

# Row 4 - 鹏华创新成长混合C
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'011461"
$q1.Range("C4").Value = "鹏华创新成长混合C"
$q1.Range("D4").Value = "'0.40"
$q1.Range("E4").Value = "'93.30"
$q1.Range("F4").Value = "'4.34"
$q1.Range("G4").Value = "'0.0174"
$q1.Range("H4").Value = 7

# Re-apply the canonical formatting (bold font + border + centred
# alignment / plain default style) that the apostrophe-forced text
# entries above leave as a bare quote-prefix style: copy it straight
# from the sibling "2021-Q4" sheet, which already carries it on the
# same cells.
$sheetQ4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$sheetQ4.Range("A2:A4").Copy()
$q1.Range("A2:A4").PasteSpecial(-4122)

$sheetQ4.Range("B2:G4").Copy()
$q1.Range("B2:G4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Brand-new "总计" sheet (cloned from "2022-Q1" so it keeps the
#    original sheet-level page setup / outline properties), appended
#    right after "2022-Q1", with the refreshed running summary:
#    2022-Q1 first, then 2021-Q4.
# ---------------------------------------------------------------------
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item(3)
$total.Name = "总计"

# Clear the cloned 2022-Q1 detail content before writing the summary
# table (the clone starts as a 4-row x 7-column holdings table; the
# summary table only needs 3 rows x 3 columns).
$total.Cells.Clear()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.61

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 9
$total.Range("D3").Value = 1.52

$sheetQ4.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$sheetQ4.Range("A2:A2").Copy()
$total.Range("A2:A3").PasteSpecial(-4122)

$sheetQ4.Range("A1").Select()
